# Remove the "신한제12호스팩" (2024-04-02 / 2024-04-05) IPO entry from each of the
# three data sheets. This row is row 21 on sheet "01_리그테이블", row 15 on sheet
# "02_통합집계_Rawdata", and row 14 on sheet "03_IPO현황_Summary".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("01_리그테이블")
$ws1.Rows.Item(21).Delete()

$ws2 = $wb.Worksheets.Item("02_통합집계_Rawdata")
$ws2.Rows.Item(15).Delete()

$ws3 = $wb.Worksheets.Item("03_IPO현황_Summary")
$ws3.Rows.Item(14).Delete()
